# Update diagrams in Developer Guide
#
# The high level sequence diagram slide renamed the "address book"
# domain concepts to "wish book" / "wish" across several shapes:
#   - deletePerson(p)                  -> deleteWish(p)
#   - post(AddressBookChangedEvent)    -> post(WishBookChangedEvent)
#   - handleAddresssBookChangedEvent() -> handleWishBookChangedEvent()

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    foreach ($shp in $slide.Shapes) {
        if ($shp.Id -eq $id) {
            return $shp
        }
    }
    return $null
}

function Replace-ShapeSubstring($shape, $oldStr, $newStr) {
    $tr = $shape.TextFrame.TextRange
    $full = $tr.Text
    $idx = $full.IndexOf($oldStr)
    if ($idx -lt 0) {
        return
    }
    $start = $idx + 1
    $len = $oldStr.Length
    $sub = $tr.Characters($start, $len)
    $sub.Text = $newStr
}

# Shape 29 "TextBox 28": deletePerson(p) -> deleteWish(p)
$shp29 = Get-ShapeById $s 29
Replace-ShapeSubstring $shp29 "deletePerson" "deleteWish"

# Shape 33 "TextBox 32": post(AddressBookChangedEvent) -> post(WishBookChangedEvent)
$shp33 = Get-ShapeById $s 33
Replace-ShapeSubstring $shp33 "AddressBookChangedEvent" "WishBookChangedEvent"

# Shape 62 "TextBox 61": post(AddressBookChangedEvent) -> post(WishBookChangedEvent)
$shp62 = Get-ShapeById $s 62
Replace-ShapeSubstring $shp62 "AddressBookChangedEvent" "WishBookChangedEvent"

# Shape 74 "TextBox 73": handleAddresssBookChangedEvent() -> handleWishBookChangedEvent()
$shp74 = Get-ShapeById $s 74
Replace-ShapeSubstring $shp74 "handleAddresssBookChangedEvent" "handleWishBookChangedEvent"

# Shape 50 "TextBox 49": handleAddresssBookChangedEvent() -> handleWishBookChangedEvent()
$shp50 = Get-ShapeById $s 50
Replace-ShapeSubstring $shp50 "handleAddresssBookChangedEvent" "handleWishBookChangedEvent"
